$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.806.62'
$ws.Range("E2").Value = '  +3.38%  '
$ws.Range("D3").Value = '3.288.54'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.53'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.97'
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").Value = '  +3.07%  '
$ws.Range("D9").Value = '3.280.81'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("E10").Value = '  +1.08%  '
$ws.Range("E11").Value = '  +1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.69'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '692.28'
$ws.Range("E14").Value = '  +12.89%  '
$ws.Range("D15").Value = '3.812.16'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.38'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '67.826.60'
$ws.Range("E17").Value = '  +3.19%  '
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").Value = '3.283.27'
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.43'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.80'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.895'
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.06'
$ws.Range("E23").Value = '  -5.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.20'
$ws.Range("E24").Value = '  +5.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.88'
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.74'
$ws.Range("E27").Value = '  +2.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.38'
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.22'
$ws.Range("E29").Value = '  +8.27%  '
$ws.Range("E30").Value = '  +1.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.74'
$ws.Range("E31").Value = '  +4.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '582.34'
$ws.Range("E32").Value = '  +7.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.88'
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D34").Value = '3.858.93'
$ws.Range("E34").Value = '  +1.81%  '
$ws.Range("E35").Value = '  +1.42%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.41'
$ws.Range("E37").Value = '  -7.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.44'
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.23'
$ws.Range("E40").Value = '  +3.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.62'
$ws.Range("E41").Value = '  +2.40%  '
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '32.18'
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0412'
$ws.Range("E46").Value = '  +1.94%  '
$ws.Range("E47").Value = '  +2.32%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.40'
$ws.Range("E48").Value = '  +10.04%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '129.12'
$ws.Range("E51").Value = '  +1.11%  '
